$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" date column (C) for all data rows (2-458)
# from serial date 45179 (2023-09-10) to 45180 (2023-09-11)
$ws.Range("C2:C458").Value = 45180
